$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove all existing hyperlinks (will re-add for the final 16 rows below)
$ws.Range("K2").Hyperlinks.Delete()

# Force column I (Exp) to text so numeric-looking values ("5","3",...) stay as text
$ws.Range("I2:I17").NumberFormat = "@"

# Row 2: Royce O'Neale
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = "Royce O'Neale"
$ws.Range("D2").Value = "SF"
$ws.Range("E2").Value = "6-4"
$ws.Range("F2").Value = 226
$ws.Range("G2").Value = "June 5, 1993"
$ws.Range("H2").Value = "us"
$ws.Range("I2").Value = "5"
$ws.Range("J2").Value = "Denver, Baylor"
$ws.Range("K2").Value = "https://www.basketball-reference.com/players/o/onealro01.html"

# Row 3: Nic Claxton
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 11
$ws.Range("C3").Value = "Nic Claxton"
$ws.Range("D3").Value = "C"
$ws.Range("E3").Value = "6-11"
$ws.Range("F3").Value = 215
$ws.Range("G3").Value = "April 17, 1999"
$ws.Range("H3").Value = "us"
$ws.Range("I3").Value = "3"
$ws.Range("J3").Value = "Georgia"
$ws.Range("K3").Value = "https://www.basketball-reference.com/players/c/claxtni01.html"

# Row 4: Joe Harris
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 12
$ws.Range("C4").Value = "Joe Harris"
$ws.Range("D4").Value = "SG"
$ws.Range("E4").Value = "6-6"
$ws.Range("F4").Value = 220
$ws.Range("G4").Value = "September 6, 1991"
$ws.Range("H4").Value = "us"
$ws.Range("I4").Value = "8"
$ws.Range("J4").Value = "Virginia"
$ws.Range("K4").Value = "https://www.basketball-reference.com/players/h/harrijo01.html"

# Row 5: Edmond Sumner
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = "Edmond Sumner"
$ws.Range("D5").Value = "SG"
$ws.Range("E5").Value = "6-4"
$ws.Range("F5").Value = 196
$ws.Range("G5").Value = "December 31, 1995"
$ws.Range("H5").Value = "us"
$ws.Range("I5").Value = "4"
$ws.Range("J5").Value = "Xavier"
$ws.Range("K5").Value = "https://www.basketball-reference.com/players/s/sumneed01.html"

# Row 6: Yuta Watanabe
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 18
$ws.Range("C6").Value = "Yuta Watanabe"
$ws.Range("D6").Value = "SF"
$ws.Range("E6").Value = "6-9"
$ws.Range("F6").Value = 215
$ws.Range("G6").Value = "October 13, 1994"
$ws.Range("H6").Value = "jp"
$ws.Range("I6").Value = "4"
$ws.Range("J6").Value = "George Washington"
$ws.Range("K6").Value = "https://www.basketball-reference.com/players/w/watanyu01.html"

# Row 7: Ben Simmons
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 10
$ws.Range("C7").Value = "Ben Simmons"
$ws.Range("D7").Value = "PG"
$ws.Range("E7").Value = "6-10"
$ws.Range("F7").Value = 240
$ws.Range("G7").Value = "July 20, 1996"
$ws.Range("H7").Value = "au"
$ws.Range("I7").Value = "4"
$ws.Range("J7").Value = "LSU"
$ws.Range("K7").Value = "https://www.basketball-reference.com/players/s/simmobe01.html"

# Row 8: Cam Thomas
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 24
$ws.Range("C8").Value = "Cam Thomas"
$ws.Range("D8").Value = "SG"
$ws.Range("E8").Value = "6-4"
$ws.Range("F8").Value = 210
$ws.Range("G8").Value = "October 13, 2001"
$ws.Range("H8").Value = "jp"
$ws.Range("I8").Value = "1"
$ws.Range("J8").Value = "LSU"
$ws.Range("K8").Value = "https://www.basketball-reference.com/players/t/thomaca02.html"

# Row 9: Seth Curry
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 30
$ws.Range("C9").Value = "Seth Curry"
$ws.Range("D9").Value = "SG"
$ws.Range("E9").Value = "6-2"
$ws.Range("F9").Value = 185
$ws.Range("G9").Value = "August 23, 1990"
$ws.Range("H9").Value = "us"
$ws.Range("I9").Value = "8"
$ws.Range("J9").Value = "Liberty, Duke"
$ws.Range("K9").Value = "https://www.basketball-reference.com/players/c/curryse01.html"

# Row 10: Patty Mills
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 8
$ws.Range("C10").Value = "Patty Mills"
$ws.Range("D10").Value = "PG"
$ws.Range("E10").Value = "6-1"
$ws.Range("F10").Value = 180
$ws.Range("G10").Value = "August 11, 1988"
$ws.Range("H10").Value = "au"
$ws.Range("I10").Value = "13"
$ws.Range("J10").Value = "Saint Mary's"
$ws.Range("K10").Value = "https://www.basketball-reference.com/players/m/millspa02.html"

# Row 11: Day'Ron Sharpe
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 20
$ws.Range("C11").Value = "Day'Ron Sharpe"
$ws.Range("D11").Value = "C"
$ws.Range("E11").Value = "6-11"
$ws.Range("F11").Value = 265
$ws.Range("G11").Value = "November 6, 2001"
$ws.Range("H11").Value = "us"
$ws.Range("I11").Value = "1"
$ws.Range("J11").Value = "UNC"
$ws.Range("K11").Value = "https://www.basketball-reference.com/players/s/sharpda01.html"

# Row 12: David Duke Jr. (TW)
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = 6
$ws.Range("C12").Value = "David Duke Jr. (TW)"
$ws.Range("D12").Value = "SG"
$ws.Range("E12").Value = "6-5"
$ws.Range("F12").Value = 205
$ws.Range("G12").Value = "October 13, 1999"
$ws.Range("H12").Value = "us"
$ws.Range("I12").Value = "1"
$ws.Range("J12").Value = "Providence"
$ws.Range("K12").Value = "https://www.basketball-reference.com/players/d/dukeda01.html"

# Row 13: Spencer Dinwiddie
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = 26
$ws.Range("C13").Value = "Spencer Dinwiddie"
$ws.Range("D13").Value = "PG"
$ws.Range("E13").Value = "6-5"
$ws.Range("F13").Value = 215
$ws.Range("G13").Value = "April 6, 1993"
$ws.Range("H13").Value = "us"
$ws.Range("I13").Value = "8"
$ws.Range("J13").Value = "Colorado"
$ws.Range("K13").Value = "https://www.basketball-reference.com/players/d/dinwisp01.html"

# Row 14: Dorian Finney-Smith
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = 28
$ws.Range("C14").Value = "Dorian Finney-Smith"
$ws.Range("D14").Value = "PF"
$ws.Range("E14").Value = "6-7"
$ws.Range("F14").Value = 220
$ws.Range("G14").Value = "May 4, 1993"
$ws.Range("H14").Value = "us"
$ws.Range("I14").Value = "6"
$ws.Range("J14").Value = "Florida"
$ws.Range("K14").Value = "https://www.basketball-reference.com/players/f/finnedo01.html"

# Row 15: Dru Smith (TW)
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = 9
$ws.Range("C15").Value = "Dru Smith (TW)"
$ws.Range("D15").Value = "SG"
$ws.Range("E15").Value = "6-3"
$ws.Range("F15").Value = 203
$ws.Range("G15").Value = "December 30, 1997"
$ws.Range("H15").Value = "us"
$ws.Range("I15").Value = "R"
$ws.Range("J15").Value = "University of Evansville, Missouri"
$ws.Range("K15").Value = "https://www.basketball-reference.com/players/s/smithdr01.html"

# Row 16: Mikal Bridges
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = "Mikal Bridges"
$ws.Range("D16").Value = "SF"
$ws.Range("E16").Value = "6-6"
$ws.Range("F16").Value = 209
$ws.Range("G16").Value = "August 30, 1996"
$ws.Range("H16").Value = "us"
$ws.Range("I16").Value = "4"
$ws.Range("J16").Value = "Villanova"
$ws.Range("K16").Value = "https://www.basketball-reference.com/players/b/bridgmi01.html"

# Row 17: Cameron Johnson
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = 2
$ws.Range("C17").Value = "Cameron Johnson"
$ws.Range("D17").Value = "PF"
$ws.Range("E17").Value = "6-8"
$ws.Range("F17").Value = 210
$ws.Range("G17").Value = "March 3, 1996"
$ws.Range("H17").Value = "us"
$ws.Range("I17").Value = "3"
$ws.Range("J17").Value = "Pitt, UNC"
$ws.Range("K17").Value = "https://www.basketball-reference.com/players/j/johnsca02.html"

# Clear the forced text format on column I so the style matches the original (no explicit numFmt)
$ws.Range("I2:I17").ClearFormats()

# Remove the now-obsolete 17th player row (Alondes Williams)
$ws.Range("A18:K18").Delete()

# Re-create the hyperlinks for the bbref url column, in row order
$ws.Hyperlinks.Add($ws.Range("K2"), "https://www.basketball-reference.com/players/o/onealro01.html")
$ws.Hyperlinks.Add($ws.Range("K3"), "https://www.basketball-reference.com/players/c/claxtni01.html")
$ws.Hyperlinks.Add($ws.Range("K4"), "https://www.basketball-reference.com/players/h/harrijo01.html")
$ws.Hyperlinks.Add($ws.Range("K5"), "https://www.basketball-reference.com/players/s/sumneed01.html")
$ws.Hyperlinks.Add($ws.Range("K6"), "https://www.basketball-reference.com/players/w/watanyu01.html")
$ws.Hyperlinks.Add($ws.Range("K7"), "https://www.basketball-reference.com/players/s/simmobe01.html")
$ws.Hyperlinks.Add($ws.Range("K8"), "https://www.basketball-reference.com/players/t/thomaca02.html")
$ws.Hyperlinks.Add($ws.Range("K9"), "https://www.basketball-reference.com/players/c/curryse01.html")
$ws.Hyperlinks.Add($ws.Range("K10"), "https://www.basketball-reference.com/players/m/millspa02.html")
$ws.Hyperlinks.Add($ws.Range("K11"), "https://www.basketball-reference.com/players/s/sharpda01.html")
$ws.Hyperlinks.Add($ws.Range("K12"), "https://www.basketball-reference.com/players/d/dukeda01.html")
$ws.Hyperlinks.Add($ws.Range("K13"), "https://www.basketball-reference.com/players/d/dinwisp01.html")
$ws.Hyperlinks.Add($ws.Range("K14"), "https://www.basketball-reference.com/players/f/finnedo01.html")
$ws.Hyperlinks.Add($ws.Range("K15"), "https://www.basketball-reference.com/players/s/smithdr01.html")
$ws.Hyperlinks.Add($ws.Range("K16"), "https://www.basketball-reference.com/players/b/bridgmi01.html")
$ws.Hyperlinks.Add($ws.Range("K17"), "https://www.basketball-reference.com/players/j/johnsca02.html")
